# Updated cryptos list on Mon Feb 27 23:24:41 UTC 2023 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns of the crypto table
# with the latest scraped figures, and swaps the Hedera / InternetComputer
# rows (39-40) back into rank order.
#
# Note: several Price values look numeric (e.g. "0.9983", "304.50") but the
# source data stores them as plain text (leading/trailing zeros, no true
# numeric semantics). A leading apostrophe forces Excel to keep them as text
# instead of silently converting them to floating point numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.485.37'
$ws.Range('E2').Value = '  -0.26%  '
$ws.Range('D3').Value = '1.631.76'
$ws.Range('E3').Value = '  -0.50%  '
$ws.Range('D4').Value = '''0.9983'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '''0.9989'
$ws.Range('D6').Value = '''304.50'
$ws.Range('E6').Value = '  -1.48%  '
$ws.Range('D7').Value = '''0.3792'
$ws.Range('E7').Value = '  +0.60%  '
$ws.Range('D8').Value = '''52.05'
$ws.Range('E8').Value = '  -1.44%  '
$ws.Range('D9').Value = '''0.3641'
$ws.Range('E9').Value = '  -1.16%  '
$ws.Range('D10').Value = '''1.235'
$ws.Range('E10').Value = '  -3.36%  '
$ws.Range('D11').Value = '''0.08118'
$ws.Range('E11').Value = '  -1.03%  '
$ws.Range('D12').Value = '''0.9983'
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('D13').Value = '''22.72'
$ws.Range('E13').Value = '  -2.06%  '
$ws.Range('D14').Value = '''6.581'
$ws.Range('E14').Value = '  -1.31%  '
$ws.Range('D15').Value = '''0.00001249'
$ws.Range('E15').Value = '  -2.46%  '
$ws.Range('D16').Value = '''7.257'
$ws.Range('E16').Value = '  -2.80%  '
$ws.Range('D17').Value = '1.622.69'
$ws.Range('E17').Value = '  -1.07%  '
$ws.Range('D18').Value = '''93.79'
$ws.Range('E18').Value = '  -1.34%  '
$ws.Range('D19').Value = '''0.06940'
$ws.Range('E19').Value = '  -0.20%  '
$ws.Range('D20').Value = '''17.95'
$ws.Range('E20').Value = '  -2.46%  '
$ws.Range('D21').Value = '''6.445'
$ws.Range('E21').Value = '  -2.08%  '
$ws.Range('D22').Value = '''0.9994'
$ws.Range('E22').Value = '  +0.16%  '
$ws.Range('D23').Value = '23.492.28'
$ws.Range('E23').Value = '  -0.26%  '
$ws.Range('D24').Value = '''12.77'
$ws.Range('E24').Value = '  -1.46%  '
$ws.Range('E25').Value = '  +3.40%  '
$ws.Range('D26').Value = '''2.417'
$ws.Range('E26').Value = '  +0.43%  '
$ws.Range('D27').Value = '''21.25'
$ws.Range('E27').Value = '  -0.70%  '
$ws.Range('D28').Value = '''149.74'
$ws.Range('E28').Value = '  -0.99%  '
$ws.Range('D29').Value = '''5.284'
$ws.Range('E29').Value = '  -0.83%  '
$ws.Range('D30').Value = '''134.93'
$ws.Range('E30').Value = '  -0.71%  '
$ws.Range('D31').Value = '''2.301'
$ws.Range('D32').Value = '1.797.96'
$ws.Range('E32').Value = '  -1.25%  '
$ws.Range('D33').Value = '''6.839'
$ws.Range('E33').Value = '  -0.47%  '
$ws.Range('E34').Value = '  +5.82%  '
$ws.Range('D35').Value = '''0.9608'
$ws.Range('E35').Value = '  -1.66%  '
$ws.Range('D36').Value = '''0.02807'
$ws.Range('E36').Value = '  -0.38%  '
$ws.Range('D37').Value = '''0.2541'
$ws.Range('E37').Value = '  -0.07%  '
$ws.Range('D38').Value = '''0.08839'
$ws.Range('E38').Value = '  -0.32%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '''0.07211'
$ws.Range('E39').Value = '  -3.49%  '
$ws.Range('B40').Value = 'InternetComputer(DFINITY)'
$ws.Range('C40').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D40').Value = '''6.124'
$ws.Range('E40').Value = '  -1.72%  '
$ws.Range('D41').Value = '''0.7111'
$ws.Range('E41').Value = '  -0.84%  '
$ws.Range('D42').Value = '''1.357'
$ws.Range('E42').Value = '  -3.16%  '
$ws.Range('D43').Value = '''16.27'
$ws.Range('E43').Value = '  +0.48%  '
$ws.Range('E44').Value = '  -2.23%  '
$ws.Range('D45').Value = '''0.6530'
$ws.Range('E45').Value = '  -1.45%  '
$ws.Range('D46').Value = '''2.341'
$ws.Range('E46').Value = '  -0.79%  '
$ws.Range('D47').Value = '''0.9984'
$ws.Range('E47').Value = '  -0.05%  '
$ws.Range('D48').Value = '''4.004'
$ws.Range('E48').Value = '  -1.12%  '
$ws.Range('D49').Value = '''0.08015'
$ws.Range('E49').Value = '  -0.50%  '
$ws.Range('E50').Value = '  -1.07%  '
$ws.Range('D51').Value = '''125.78'
$ws.Range('E51').Value = '  -4.01%  '
